$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data set (Vencimento, Taxa, Data de Salvamento) - 42 data rows starting at row 2
$data = @(
    @(45792, 0, "2025-04-04 13:12:50"),
    @(45792, 1, "2025-04-04 13:13:02"),
    @(45792, 3, "2025-04-04 13:13:12"),
    @(46249, 0, "2025-04-04 13:12:50"),
    @(46249, 1, "2025-04-04 13:13:02"),
    @(46249, 3, "2025-04-04 13:13:12"),
    @(46522, 3, "2025-04-04 13:13:12"),
    @(46522, 0, "2025-04-04 13:12:50"),
    @(46522, 1, "2025-04-04 13:13:02"),
    @(46980, 1, "2025-04-04 13:13:02"),
    @(46980, 3, "2025-04-04 13:13:12"),
    @(46980, 0, "2025-04-04 13:12:50"),
    @(47253, 0, "2025-04-04 13:12:50"),
    @(47253, 1, "2025-04-04 13:13:02"),
    @(47253, 3, "2025-04-04 13:13:12"),
    @(47710, 0, "2025-04-04 13:12:50"),
    @(47710, 1, "2025-04-04 13:13:02"),
    @(47710, 3, "2025-04-04 13:13:12"),
    @(48441, 0, "2025-04-04 13:12:50"),
    @(48441, 1, "2025-04-04 13:13:02"),
    @(48441, 3, "2025-04-04 13:13:12"),
    @(48714, 3, "2025-04-04 13:13:12"),
    @(48714, 0, "2025-04-04 13:12:50"),
    @(48714, 1, "2025-04-04 13:13:02"),
    @(49444, 1, "2025-04-04 13:13:02"),
    @(49444, 0, "2025-04-04 13:12:50"),
    @(49444, 3, "2025-04-04 13:13:12"),
    @(51363, 3, "2025-04-04 13:13:12"),
    @(51363, 0, "2025-04-04 13:12:50"),
    @(51363, 1, "2025-04-04 13:13:02"),
    @(53097, 3, "2025-04-04 13:13:12"),
    @(53097, 1, "2025-04-04 13:13:02"),
    @(53097, 0, "2025-04-04 13:12:50"),
    @(55015, 0, "2025-04-04 13:12:50"),
    @(55015, 1, "2025-04-04 13:13:02"),
    @(55015, 3, "2025-04-04 13:13:12"),
    @(56749, 0, "2025-04-04 13:12:50"),
    @(56749, 1, "2025-04-04 13:13:02"),
    @(56749, 3, "2025-04-04 13:13:12"),
    @(58668, 0, "2025-04-04 13:12:50"),
    @(58668, 1, "2025-04-04 13:13:02"),
    @(58668, 3, "2025-04-04 13:13:12")
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $row = $row + 1
}
